$d = $word.ActiveDocument

$d.Content.Find.Execute("13×30=", $true, $false, $false, $false, $false, $true, 1, $false, "40×37=", 2) | Out-Null
$d.Content.Find.Execute("80×91=", $true, $false, $false, $false, $false, $true, 1, $false, "88×42=", 2) | Out-Null
$d.Content.Find.Execute("68×37=", $true, $false, $false, $false, $false, $true, 1, $false, "97×88=", 2) | Out-Null
$d.Content.Find.Execute("54×77=", $true, $false, $false, $false, $false, $true, 1, $false, "89×28=", 2) | Out-Null
$d.Content.Find.Execute("30×51=", $true, $false, $false, $false, $false, $true, 1, $false, "88×91=", 2) | Out-Null
$d.Content.Find.Execute("78×55=", $true, $false, $false, $false, $false, $true, 1, $false, "42×95=", 2) | Out-Null
$d.Content.Find.Execute("11×41=", $true, $false, $false, $false, $false, $true, 1, $false, "48×89=", 2) | Out-Null
$d.Content.Find.Execute("18×37=", $true, $false, $false, $false, $false, $true, 1, $false, "70×66=", 2) | Out-Null
$d.Content.Find.Execute("40×21=", $true, $false, $false, $false, $false, $true, 1, $false, "18×95=", 2) | Out-Null
$d.Content.Find.Execute("77×41=", $true, $false, $false, $false, $false, $true, 1, $false, "31×78=", 2) | Out-Null
$d.Content.Find.Execute("78×23=", $true, $false, $false, $false, $false, $true, 1, $false, "52×32=", 2) | Out-Null
$d.Content.Find.Execute("29×96=", $true, $false, $false, $false, $false, $true, 1, $false, "52×19=", 2) | Out-Null
$d.Content.Find.Execute("69×60=", $true, $false, $false, $false, $false, $true, 1, $false, "28×92=", 2) | Out-Null
$d.Content.Find.Execute("46×53=", $true, $false, $false, $false, $false, $true, 1, $false, "49×58=", 2) | Out-Null
$d.Content.Find.Execute("52×87=", $true, $false, $false, $false, $false, $true, 1, $false, "33×74=", 2) | Out-Null
$d.Content.Find.Execute("95×46=", $true, $false, $false, $false, $false, $true, 1, $false, "37×63=", 2) | Out-Null
$d.Content.Find.Execute("93×83=", $true, $false, $false, $false, $false, $true, 1, $false, "46×81=", 2) | Out-Null
$d.Content.Find.Execute("38×66=", $true, $false, $false, $false, $false, $true, 1, $false, "52×82=", 2) | Out-Null
$d.Content.Find.Execute("87×85=", $true, $false, $false, $false, $false, $true, 1, $false, "38×95=", 2) | Out-Null
$d.Content.Find.Execute("39×22=", $true, $false, $false, $false, $false, $true, 1, $false, "88×44=", 2) | Out-Null
$d.Content.Find.Execute("99×59=", $true, $false, $false, $false, $false, $true, 1, $false, "94×41=", 2) | Out-Null
$d.Content.Find.Execute("82×18=", $true, $false, $false, $false, $false, $true, 1, $false, "27×14=", 2) | Out-Null
$d.Content.Find.Execute("76×86=", $true, $false, $false, $false, $false, $true, 1, $false, "66×22=", 2) | Out-Null
$d.Content.Find.Execute("26×44=", $true, $false, $false, $false, $false, $true, 1, $false, "29×71=", 2) | Out-Null
$d.Content.Find.Execute("66×47=", $true, $false, $false, $false, $false, $true, 1, $false, "74×13=", 2) | Out-Null
